# Apply updated crypto price/volume figures to Sheet1 (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.132.20"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "2.026.09"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.48%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0792"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("D12").Value = "2.322.83"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.745"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").Value = "2.018.41"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "36.981.15"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("E25").Value = "  -6.10%  "
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("E28").Value = "  -4.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("E30").Value = "  -3.80%  "
$ws.Range("E31").Value = "  -4.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "1.506.85"
$ws.Range("E40").Value = "  +3.79%  "
$ws.Range("E41").Value = "  -7.34%  "
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0932"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.73%  "
$ws.Range("E46").Value = "  -5.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.56%  "
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").Value = "2.211.00"
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.01%  "
